$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("ECs", "Cthrc1", "Fzd3", "ECs", 3, 1, 0.6191033333333333, 1.85731, 0.007929637811697733, 0.007962490129789305, 3, 1, 0.2078313333333333, 0.623494, 0.08621557350328635, 0.112461889302165, 0.1286690712377778, 1.15802164114, 0.0006836582716088646, 0.0008954766835459462),
    @("ECs", "Cthrc1", "Fzd3", "FAPs", 3, 1, 0.6191033333333333, 1.85731, 0.007929637811697733, 0.007962490129789305, 3, 1, 0.498127, 1.494381, 0.206640184103479, 0.2695469573039334, 0.3083920861233334, 2.77552877511, 0.001638581817283128, 0.002146264987047309),
    @("ECs", "Cthrc1", "Fzd3", "M2", 3, 1, 0.6191033333333333, 1.85731, 0.007929637811697733, 0.007962490129789305, 1, 0.3333333333333333, 0.016887, 0.050661, 0.007005307459654767, 0.009137909545139137, 0.01045479799, 0.09409318190999999, 0.00005554955091464663, 0.00007276051456007786),
    @("ECs", "Cthrc1", "Fzd3", "sCs", 3, 1, 0.6191033333333333, 1.85731, 0.007929637811697733, 0.007962490129789305, 2, 1, 1.6877555, 3.375511, 0.7001389349335798, 0.6088532438487625, 1.044895055901667, 6.269370335410001, 0.005551848171891093, 0.004847987944635973),
    @("FAPs", "Cthrc1", "Fzd3", "ECs", 3, 1, 76.48912033333333, 229.467361, 0.9796927073757713, 0.9837515530909214, 3, 1, 0.2078313333333333, 0.623494, 0.08621557350328635, 0.112461889302165, 15.89683586437044, 143.071522779334, 0.08446476862338942, 0.1106345582645441),
    @("FAPs", "Cthrc1", "Fzd3", "FAPs", 3, 1, 76.48912033333333, 229.467361, 0.9796927073757713, 0.9837515530909214, 3, 1, 0.498127, 1.494381, 0.206640184103479, 0.2695469573039334, 38.10129604428234, 342.911664398541, 0.2024438814169652, 0.2651672378786768),
    @("FAPs", "Cthrc1", "Fzd3", "M2", 3, 1, 76.48912033333333, 229.467361, 0.9796927073757713, 0.9837515530909214, 1, 0.3333333333333333, 0.016887, 0.050661, 0.007005307459654767, 0.009137909545139137, 1.291671775069, 11.625045975621, 0.006863048631148866, 0.008989432707034981),
    @("FAPs", "Cthrc1", "Fzd3", "sCs", 3, 1, 76.48912033333333, 229.467361, 0.9796927073757713, 0.9837515530909214, 2, 1, 1.6877555, 3.375511, 0.7001389349335798, 0.6088532438487625, 129.0949335327452, 774.569601196471, 0.6859210087042678, 0.5989603242406656),
    @("sCs", "Cthrc1", "Fzd3", "ECs", 2, 1, 0.9663805, 1.932761, 0.01237765481253107, 0.008285956779289245, 3, 1, 0.2078313333333333, 0.623494, 0.08621557350328635, 0.112461889302165, 0.2008441478223333, 1.205064886934, 0.001067146608288078, 0.0009318543540749505),
    @("sCs", "Cthrc1", "Fzd3", "FAPs", 2, 1, 0.9663805, 1.932761, 0.01237765481253107, 0.008285956779289245, 3, 1, 0.498127, 1.494381, 0.206640184103479, 0.2695469573039334, 0.4813802193235, 2.888281315941, 0.002557720869230734, 0.002233454438209316),
    @("sCs", "Cthrc1", "Fzd3", "M2", 2, 1, 0.9663805, 1.932761, 0.01237765481253107, 0.008285956779289245, 1, 0.3333333333333333, 0.016887, 0.050661, 0.007005307459654767, 0.009137909545139137, 0.0163192675035, 0.097915605021, 0.00008670927759125563, 0.00007571632354407753),
    @("sCs", "Cthrc1", "Fzd3", "sCs", 2, 1, 0.9663805, 1.932761, 0.01237765481253107, 0.008285956779289245, 2, 1, 1.6877555, 3.375511, 0.7001389349335798, 0.6088532438487625, 1.63101400396775, 6.524056015871, 0.008666078057421003, 0.005044931663460902)
)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($startRow + $r, $c + 1).Value = $row[$c]
    }
}
